# Update the OncoFAIR MA Element Traceability metadata sheet:
#  1. Bump the "Date" property value (row 8) to the new generation timestamp.
#  2. Insert a new "Jurisdiction" property row (with an empty value) right
#     after "Contact" (row 10) and before "Description" (shifts everything
#     below down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update the Date value cell.
$ws.Cells.Item(8, 2).Value = "2024-07-01T07:50:29+00:00"

# 2. Insert a new row above "Description" (row 11) so the new row becomes
#    row 11, pushing Description/Purpose/Copyright/FHIR Version/... down by one.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row with the Jurisdiction property (empty value).
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

# Match the formatting used by every other data row (the inserted row otherwise
# comes in with the default style instead of the bordered "data row" style).
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
